# Gravity Boii Runner - Estimation workbook update
# Reflects: "Fixed Collision of powerups, tapping functionality. Gun functionality left"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (BGR-encoded OLE values matching the workbook's existing fills)
$GREEN  = 5296274   # FF92D050 - "Done"
$YELLOW = 65535      # FFFFFF00 - "WIP" (unchanged, kept for reference)

# ---------------------------------------------------------------------------
# Touch Controller section (rows 17-20)
# Swiping for movements (18) and Left/right for shifting (19) are now DONE
# -> fill green. Up/down for dashing (20) is still outstanding -> stays yellow.
# ---------------------------------------------------------------------------
$ws.Range("B18:D19").Interior.Color = $GREEN

# Time taken updates: dashing input work took longer, shifting work took less
$ws.Range("D19").Value2 = 1
$ws.Range("D20").Value2 = 2

# Updated remark for the Touch Controller group: still need gun pickup code
$ws.Range("E17").Value2 = "Create the code for gun pickup"

# ---------------------------------------------------------------------------
# Power ups section (rows 29-32) - collision issue fixed, section now Done
# ---------------------------------------------------------------------------
$ws.Range("A29").Interior.Color = $GREEN
$ws.Range("C29:D29").Interior.Color = $GREEN
$ws.Range("E29").Interior.Color = $GREEN

$ws.Range("B30:D30").Interior.Color = $GREEN
$ws.Range("E30").Interior.Color = $GREEN

$ws.Range("B31:D31").Interior.Color = $GREEN
$ws.Range("E31").Interior.Color = $GREEN

$ws.Range("B32:D32").Interior.Color = $GREEN
$ws.Range("E32").Interior.Color = $GREEN

# Gun pickup actual time taken so far
$ws.Range("C32").Value2 = 1
$ws.Range("D32").Value2 = 0.5

# Mark the Power ups group as Done
$ws.Range("E29").Value2 = "Done"

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("C23").Select()
